$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume) hold text-formatted numeric-looking
# strings (e.g. "1.002", "27.815.71") and percentages. Force the whole
# data range to Text first so Excel does not reinterpret values that
# look like numbers, then restore the default (unstyled) appearance.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = "27.815.71"
$ws.Range("E2").Value = "  -0.46%  "
$ws.Range("D3").Value = "1.908.34"
$ws.Range("E3").Value = "  +0.29%  "
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.38%  "
$ws.Range("D5").Value = "313.13"
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  -0.31%  "
$ws.Range("D7").Value = "0.4995"
$ws.Range("E7").Value = "  +3.66%  "
$ws.Range("D8").Value = "0.3783"
$ws.Range("E8").Value = "  -0.23%  "
$ws.Range("D9").Value = "0.07270"
$ws.Range("E9").Value = "  -1.28%  "
$ws.Range("D10").Value = "21.17"
$ws.Range("E10").Value = "  +2.05%  "
$ws.Range("D11").Value = "0.9020"
$ws.Range("E11").Value = "  -3.16%  "
$ws.Range("D12").Value = "0.07634"
$ws.Range("E12").Value = "  -1.33%  "
$ws.Range("D13").Value = "1.895.76"
$ws.Range("E13").Value = "  -2.75%  "
$ws.Range("D14").Value = "5.469"
$ws.Range("E14").Value = "  -0.15%  "
$ws.Range("D15").Value = "92.05"
$ws.Range("E15").Value = "  +0.37%  "
$ws.Range("E16").Value = "  -0.37%  "
$ws.Range("D17").Value = "0.000008716"
$ws.Range("E17").Value = "  -1.63%  "
$ws.Range("D18").Value = "1.001"
$ws.Range("E18").Value = "  -0.18%  "
$ws.Range("D19").Value = "27.866.57"
$ws.Range("E19").Value = "  -0.48%  "
$ws.Range("D20").Value = "14.57"
$ws.Range("E20").Value = "  -0.60%  "
$ws.Range("D21").Value = "5.166"
$ws.Range("E21").Value = "  +0.35%  "
$ws.Range("D22").Value = "2.125.67"
$ws.Range("E22").Value = "  -2.85%  "
$ws.Range("E23").Value = "  -0.52%  "
$ws.Range("D24").Value = "6.595"
$ws.Range("E24").Value = "  -0.51%  "
$ws.Range("D25").Value = "153.20"
$ws.Range("E25").Value = "  -1.86%  "
$ws.Range("D26").Value = "1.849"
$ws.Range("E26").Value = "  -3.15%  "
$ws.Range("D27").Value = "2.215"
$ws.Range("E27").Value = "  +4.53%  "
$ws.Range("D28").Value = "18.36"
$ws.Range("E28").Value = "  -0.46%  "
$ws.Range("D29").Value = "115.17"
$ws.Range("E29").Value = "  -1.60%  "
$ws.Range("D30").Value = "4.874"
$ws.Range("E30").Value = "  -1.84%  "
$ws.Range("D31").Value = "0.08966"
$ws.Range("D32").Value = "3.195"
$ws.Range("D33").Value = "1.239"
$ws.Range("E33").Value = "  -0.90%  "
$ws.Range("D34").Value = "0.7882"
$ws.Range("E34").Value = "  +2.63%  "
$ws.Range("D35").Value = "4.793"
$ws.Range("E35").Value = "  +2.75%  "
$ws.Range("D36").Value = "2.642"
$ws.Range("E36").Value = "  +2.33%  "
$ws.Range("E37").Value = "  +1.21%  "
$ws.Range("D38").Value = "3.062"
$ws.Range("E38").Value = "  +2.00%  "
$ws.Range("D39").Value = "1.092"
$ws.Range("E39").Value = "  -0.98%  "
$ws.Range("E40").Value = "  +0.56%  "
$ws.Range("D41").Value = "0.05301"
$ws.Range("E41").Value = "  +0.59%  "
$ws.Range("D42").Value = "6.752"
$ws.Range("E42").Value = "  -2.52%  "
$ws.Range("D43").Value = "114.27"
$ws.Range("E43").Value = "  +3.82%  "
$ws.Range("D44").Value = "8.471"
$ws.Range("E44").Value = "  -0.12%  "
$ws.Range("D45").Value = "0.1511"
$ws.Range("E45").Value = "  -0.85%  "
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").Value = "0.4791"
$ws.Range("E46").Value = "  -0.37%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "10.50"
$ws.Range("E47").Value = "  -2.15%  "
$ws.Range("D48").Value = "1.001"
$ws.Range("E48").Value = "  -0.28%  "
$ws.Range("E49").Value = "  -0.62%  "
$ws.Range("D50").Value = "67.26"
$ws.Range("E50").Value = "  -1.01%  "
$ws.Range("E51").Value = "  -0.74%  "

$dataRange.Style = "Normal"
